# Auto-generated edit script: update cryptocurrency price/volume table
# to match the Sat Jun  1 22:55:26 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.793.69'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '3.825.45'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'602.89"
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').Value = '4.461.20'
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').Value = '3.840.98'
$ws.Range('E15').Value = '  +2.14%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = "'18.48"
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '67.834.42'
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('E18').Value = '  +1.38%  '
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('D20').Value = "'463.14"
$ws.Range('E20').Value = '  +1.51%  '
$ws.Range('D21').Value = "'9.93"
$ws.Range('E21').Value = '  -1.31%  '
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('E23').Value = '  -3.51%  '
$ws.Range('D24').Value = "'83.45"
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').Value = "'12.10"
$ws.Range('E25').Value = '  +1.61%  '
$ws.Range('E26').Value = '  -0.93%  '
$ws.Range('D27').Value = "'10.10"
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').Value = '3.974.52'
$ws.Range('E29').Value = '  +1.20%  '
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('E31').Value = '  +1.80%  '
$ws.Range('E32').Value = '  +1.99%  '
$ws.Range('D33').Value = "'29.70"
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').Value = '3.765.72'
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('D36').Value = "'9.08"
$ws.Range('E36').Value = '  -1.35%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').Value = "'3.35"
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('E41').Value = '  +1.07%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = "'29.31"
$ws.Range('E44').Value = '  +14.41%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = "'48.09"
$ws.Range('E45').Value = '  +1.94%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').Value = "'43.06"
$ws.Range('E47').Value = '  -4.97%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = "'1.41"
$ws.Range('E48').Value = '  +12.46%  '
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('D50').Value = "'148.24"
$ws.Range('E50').Value = '  -0.09%  '
